$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-3 and add new rows 4-7 per the NATMI re-analysis
# (adds the "ECs" cluster alongside existing sCs/FAPs clusters for Bdnf-Ntrk2 pairs)

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Bdnf"
$ws.Range("C2").Value = "Ntrk2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03885866666666667
$ws.Range("H2").Value = 0.116576
$ws.Range("I2").Value = 0.01924839521029073
$ws.Range("J2").Value = 0.01924839521029073
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.5561916666666666
$ws.Range("N2").Value = 1.668575
$ws.Range("O2").Value = 0.01178629317764927
$ws.Range("P2").Value = 0.01178629317764927
$ws.Range("Q2").Value = 0.02161286657777777
$ws.Range("R2").Value = 0.1945157992
$ws.Range("S2").Value = 0.0002268672291477465
$ws.Range("T2").Value = 0.0002268672291477465

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Bdnf"
$ws.Range("C3").Value = "Ntrk2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03885866666666667
$ws.Range("H3").Value = 0.116576
$ws.Range("I3").Value = 0.01924839521029073
$ws.Range("J3").Value = 0.01924839521029073
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 43.12631833333334
$ws.Range("N3").Value = 129.378955
$ws.Range("O3").Value = 0.9138925697963186
$ws.Range("P3").Value = 0.9138925697963186
$ws.Range("Q3").Value = 1.675831228675556
$ws.Range("R3").Value = 15.08248105808
$ws.Range("S3").Value = 0.01759096536318775
$ws.Range("T3").Value = 0.01759096536318775

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Bdnf"
$ws.Range("C4").Value = "Ntrk2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.03885866666666667
$ws.Range("H4").Value = 0.116576
$ws.Range("I4").Value = 0.01924839521029073
$ws.Range("J4").Value = 0.01924839521029073
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.507192333333334
$ws.Range("N4").Value = 10.521577
$ws.Range("O4").Value = 0.0743211370260321
$ws.Range("P4").Value = 0.07432113702603209
$ws.Range("Q4").Value = 0.1362848178168889
$ws.Range("R4").Value = 1.226563360352
$ws.Range("S4").Value = 0.001430562617955237
$ws.Range("T4").Value = 0.001430562617955237

# Row 5
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Bdnf"
$ws.Range("C5").Value = "Ntrk2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.979941666666667
$ws.Range("H5").Value = 5.939825
$ws.Range("I5").Value = 0.9807516047897092
$ws.Range("J5").Value = 0.9807516047897092
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.5561916666666666
$ws.Range("N5").Value = 1.668575
$ws.Range("O5").Value = 0.01178629317764927
$ws.Range("P5").Value = 0.01178629317764927
$ws.Range("Q5").Value = 1.101227055486111
$ws.Range("R5").Value = 9.911043499374999
$ws.Range("S5").Value = 0.01155942594850153
$ws.Range("T5").Value = 0.01155942594850152

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Bdnf"
$ws.Range("C6").Value = "Ntrk2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.979941666666667
$ws.Range("H6").Value = 5.939825
$ws.Range("I6").Value = 0.9807516047897092
$ws.Range("J6").Value = 0.9807516047897092
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 43.12631833333334
$ws.Range("N6").Value = 129.378955
$ws.Range("O6").Value = 0.9138925697963186
$ws.Range("P6").Value = 0.9138925697963186
$ws.Range("Q6").Value = 85.38759459809724
$ws.Range("R6").Value = 768.4883513828751
$ws.Range("S6").Value = 0.8963016044331308
$ws.Range("T6").Value = 0.8963016044331308

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Bdnf"
$ws.Range("C7").Value = "Ntrk2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.979941666666667
$ws.Range("H7").Value = 5.939825
$ws.Range("I7").Value = 0.9807516047897092
$ws.Range("J7").Value = 0.9807516047897092
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.507192333333334
$ws.Range("N7").Value = 10.521577
$ws.Range("O7").Value = 0.0743211370260321
$ws.Range("P7").Value = 0.07432113702603209
$ws.Range("Q7").Value = 6.944036233780556
$ws.Range("R7").Value = 62.496326104025
$ws.Range("S7").Value = 0.07289057440807686
$ws.Range("T7").Value = 0.07289057440807685
